# Update the "取得日時" (retrieved timestamp) column on the "ランサーズ" sheet
# for all existing data rows (2-24) to the new run timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-10-22 18:34:56"

for ($row = 2; $row -le 24; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
